$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to remain Text so the stored cell type/content matches the source data.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '26.972.26'
$ws.Range('E2').Value = '  -1.44%  '
$ws.Range('D3').Value = '1.818.60'
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = '309.06'
$ws.Range('E5').Value = '  -1.71%  '
$ws.Range('B6').Value = 'USDC'
$ws.Range('C6').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D6').Value = '1.008'
$ws.Range('E6').Value = '  -0.19%  '
$ws.Range('D7').Value = '0.4683'
$ws.Range('E7').Value = '  -1.23%  '
$ws.Range('D8').Value = '0.3650'
$ws.Range('E8').Value = '  -1.18%  '
$ws.Range('D9').Value = '0.07221'
$ws.Range('E9').Value = '  -3.28%  '
$ws.Range('D10').Value = '0.8595'
$ws.Range('E10').Value = '  -3.01%  '
$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D11').Value = '1.906.70'
$ws.Range('E11').Value = '  +1.72%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').Value = '19.75'
$ws.Range('E12').Value = '  -3.42%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = '0.07558'
$ws.Range('E13').Value = '  +3.03%  '
$ws.Range('E14').Value = '  -2.41%  '
$ws.Range('D15').Value = '91.74'
$ws.Range('E15').Value = '  -1.42%  '
$ws.Range('D16').Value = '6.467'
$ws.Range('E16').Value = '  -1.85%  '
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').Value = '0.000008619'
$ws.Range('E18').Value = '  -2.23%  '
$ws.Range('E19').Value = '  -0.33%  '
$ws.Range('E20').Value = '  -2.34%  '
$ws.Range('D21').Value = '26.601.41'
$ws.Range('E21').Value = '  -3.46%  '
$ws.Range('D22').Value = '5.133'
$ws.Range('E22').Value = '  -3.30%  '
$ws.Range('D23').Value = '10.51'
$ws.Range('E23').Value = '  -1.57%  '
$ws.Range('D24').Value = '2.066.79'
$ws.Range('E24').Value = '  -1.35%  '
$ws.Range('D26').Value = '1.852'
$ws.Range('E26').Value = '  -2.93%  '
$ws.Range('D27').Value = '18.09'
$ws.Range('D28').Value = '2.063'
$ws.Range('E28').Value = '  -3.76%  '
$ws.Range('D29').Value = '5.107'
$ws.Range('E29').Value = '  -2.56%  '
$ws.Range('D30').Value = '115.25'
$ws.Range('E30').Value = '  -1.83%  '
$ws.Range('D31').Value = '0.08877'
$ws.Range('E31').Value = '  -1.41%  '
$ws.Range('D32').Value = '2.970'
$ws.Range('E32').Value = '  +0.89%  '
$ws.Range('D33').Value = '4.412'
$ws.Range('E33').Value = '  -2.95%  '
$ws.Range('E34').Value = '  -4.37%  '
$ws.Range('D35').Value = '0.7155'
$ws.Range('E35').Value = '  -5.46%  '
$ws.Range('D36').Value = '1.077'
$ws.Range('E36').Value = '  -2.44%  '
$ws.Range('D37').Value = '0.05248'
$ws.Range('E37').Value = '  -2.01%  '
$ws.Range('E38').Value = '  -1.80%  '
$ws.Range('D39').Value = '2.917'
$ws.Range('E39').Value = '  -2.16%  '
$ws.Range('D40').Value = '2.375'
$ws.Range('E40').Value = '  -0.97%  '
$ws.Range('D41').Value = '7.124'
$ws.Range('E41').Value = '  -2.62%  '
$ws.Range('D42').Value = '0.5146'
$ws.Range('E42').Value = '  -3.29%  '
$ws.Range('B43').Value = 'Frax'
$ws.Range('C43').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D43').Value = '0.8604'
$ws.Range('E43').Value = '  -14.84%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').Value = '0.1624'
$ws.Range('E44').Value = '  -2.25%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').Value = '8.146'
$ws.Range('E45').Value = '  -4.10%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '0.4808'
$ws.Range('E46').Value = '  -2.14%  '
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').Value = '1.009'
$ws.Range('E47').Value = '  -0.23%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '10.12'
$ws.Range('E48').Value = '  -3.81%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = '102.67'
$ws.Range('E49').Value = '  -2.21%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.06255'
$ws.Range('E50').Value = '  -0.63%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = '1.616'
$ws.Range('E51').Value = '  -3.52%  '
